# Translate the "Dictionnaire de données" sheet from French field/table
# names to English ones (commit: "dictionnaire de données EN / modification
# en anglais").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table "Client" (rows 2-8) ---
$ws.Cells.Item(3, 2).Value = "last_name"
$ws.Cells.Item(4, 2).Value = "first_name"
$ws.Cells.Item(5, 2).Value = "phone"
$ws.Cells.Item(6, 2).Value = "e-mail"
$ws.Cells.Item(7, 2).Value = "adress_id"
$ws.Cells.Item(8, 2).Value = "renting_id"

# --- Table "Adresse" -> "Adress" (rows 12-18) ---
$ws.Cells.Item(12, 1).Value = "Adress"
$ws.Cells.Item(13, 2).Value = "number"
$ws.Cells.Item(14, 2).Value = "street"
$ws.Cells.Item(15, 2).Value = "postal_code"
$ws.Cells.Item(16, 2).Value = "city"
$ws.Cells.Item(17, 2).Value = "country"
$ws.Cells.Item(18, 2).Value = "complement"

# --- Table "Location" -> "Renting" (rows 22-26) ---
$ws.Cells.Item(22, 1).Value = "Renting"
$ws.Cells.Item(23, 2).Value = "start_date"
$ws.Cells.Item(24, 2).Value = "end_date"
$ws.Cells.Item(25, 2).Value = "return"
$ws.Cells.Item(26, 2).Value = "availability"

# --- Table "Auteur" -> "Author" (rows 30-34) ---
$ws.Cells.Item(30, 1).Value = "Author"
$ws.Cells.Item(31, 2).Value = "last_name"
$ws.Cells.Item(32, 2).Value = "first_name"
$ws.Cells.Item(33, 2).Value = "date_of_birth"
$ws.Cells.Item(34, 2).Value = "nationality"

# --- Table "Livre" -> "Book" (rows 38-48) ---
$ws.Cells.Item(38, 1).Value = "Book"
$ws.Cells.Item(40, 2).Value = "title"
$ws.Cells.Item(41, 2).Value = "publication_date"
$ws.Cells.Item(42, 2).Value = "number_of_page"
$ws.Cells.Item(43, 2).Value = "purchase_date"
$ws.Cells.Item(44, 2).Value = "language"
$ws.Cells.Item(45, 2).Value = "author_id"
$ws.Cells.Item(47, 2).Value = "renting_id"
$ws.Cells.Item(48, 2).Value = "copy_id"

# --- Table "Exemplaire" -> "Copy" (rows 52-53) ---
$ws.Cells.Item(52, 1).Value = "Copy"
$ws.Cells.Item(53, 2).Value = "copy_number"

# --- Table "Type" (rows 57-58) ---
$ws.Cells.Item(58, 2).Value = "type_name"

# Restore selection to the "Livre"/"Book" block and scroll so row 34 is
# the top-left visible row, matching the saved view state.
$ws.Range("A38:A48").Select()
$excel.ActiveWindow.ScrollRow = 34
